# Update the BOM: capacitor footprint/part-number change, and move the
# active-cell selection to reflect the adjusted connector position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bom")

# Rows 2-23 correspond to C1..C22 (the 1uF capacitors). Their footprint
# changed from C0402 to C0603, and the JLCPCB part number changed from
# C52923 to C15849.
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 4).Value = "C15849"
    $ws.Cells.Item($r, 3).Value = "C0603"
}

# Adjusted connector position: move the active selection to C25.
$ws.Range("C25").Select()

$wb.Save()
